$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'310.34"
$ws.Range("E2").Value = "'0.26%"
$ws.Range("G2").Value = "'3"
$ws.Range("D3").Value = "'39.86"
$ws.Range("E3").Value = "'-3.04%"
$ws.Range("G3").Value = "'3"
$ws.Range("D4").Value = "'5.092"
$ws.Range("E4").Value = "'-2.90%"
$ws.Range("G4").Value = "'3"
$ws.Range("D5").Value = "'0.07556"
$ws.Range("E5").Value = "'-1.48%"
$ws.Range("G5").Value = "'3"
$ws.Range("D6").Value = "'4.306"
$ws.Range("E6").Value = "'-0.55%"
$ws.Range("G6").Value = "'3"
$ws.Range("D7").Value = "'1.656"
$ws.Range("E7").Value = "'1.88%"
$ws.Range("G7").Value = "'3"
$ws.Range("D8").Value = "'0.9304"
$ws.Range("E8").Value = "'1.20%"
$ws.Range("G8").Value = "'3"
$ws.Range("E9").Value = "'-1.35%"
$ws.Range("G9").Value = "'3"
$ws.Range("D10").Value = "'0.1225"
$ws.Range("E10").Value = "'-2.50%"
$ws.Range("G10").Value = "'3"
$ws.Range("D11").Value = "'0.1801"
$ws.Range("E11").Value = "'-2.17%"
$ws.Range("G11").Value = "'3"
$ws.Range("D12").Value = "'0.09051"
$ws.Range("E12").Value = "'-1.38%"
$ws.Range("G12").Value = "'3"
$ws.Range("D13").Value = "'0.04146"
$ws.Range("E13").Value = "'-3.41%"
$ws.Range("G13").Value = "'3"
$ws.Range("E14").Value = "'0.19%"
$ws.Range("G14").Value = "'3"
$ws.Range("D15").Value = "'0.001281"
$ws.Range("E15").Value = "'1.66%"
$ws.Range("G15").Value = "'3"
$ws.Range("D16").Value = "'0.005769"
$ws.Range("E16").Value = "'-0.35%"
$ws.Range("G16").Value = "'3"
$ws.Range("G17").Value = "'3"
$ws.Range("E18").Value = "'-0.22%"
$ws.Range("G18").Value = "'3"
$ws.Range("D19").Value = "'0.3354"
$ws.Range("E19").Value = "'0.56%"
$ws.Range("G19").Value = "'3"
$ws.Range("D20").Value = "'7.698"
$ws.Range("E20").Value = "'6.62%"
$ws.Range("G20").Value = "'3"
$ws.Range("D21").Value = "'0.1353"
$ws.Range("E21").Value = "'-2.19%"
$ws.Range("G21").Value = "'3"
$ws.Range("D22").Value = "'0.2932"
$ws.Range("E22").Value = "'0.25%"
$ws.Range("G22").Value = "'3"
$ws.Range("D23").Value = "'0.04027"
$ws.Range("E23").Value = "'-0.96%"
$ws.Range("G23").Value = "'3"
$ws.Range("D24").Value = "'0.001267"
$ws.Range("E24").Value = "'0.48%"
$ws.Range("G24").Value = "'3"
$ws.Range("D25").Value = "'0.004084"
$ws.Range("E25").Value = "'-2.01%"
$ws.Range("G25").Value = "'3"
$ws.Range("D26").Value = "'0.0001273"
$ws.Range("E26").Value = "'0.04%"
$ws.Range("G26").Value = "'3"
$ws.Range("G27").Value = "'3"
$ws.Range("G28").Value = "'3"
$ws.Range("G29").Value = "'3"
$ws.Range("G30").Value = "'3"
$ws.Range("G31").Value = "'3"
$ws.Range("G32").Value = "'3"
$ws.Range("G33").Value = "'3"
$ws.Range("G34").Value = "'3"
$ws.Range("G35").Value = "'3"
$ws.Range("G36").Value = "'3"
$ws.Range("G37").Value = "'3"
$ws.Range("D38").Value = "'0.02427"
$ws.Range("E38").Value = "'-0.99%"
$ws.Range("G38").Value = "'3"
$ws.Range("D39").Value = "'0.05140"
$ws.Range("E39").Value = "'-3.07%"
$ws.Range("G39").Value = "'3"
$ws.Range("D40").Value = "'0.007725"
$ws.Range("E40").Value = "'-1.48%"
$ws.Range("G40").Value = "'3"
$ws.Range("D41").Value = "'0.1297"
$ws.Range("E41").Value = "'-1.35%"
$ws.Range("G41").Value = "'3"
$ws.Range("D42").Value = "'0.007676"
$ws.Range("E42").Value = "'12.60%"
$ws.Range("G42").Value = "'3"
$ws.Range("E43").Value = "'14.28%"
$ws.Range("G43").Value = "'3"
$ws.Range("D44").Value = "'0.007991"
$ws.Range("E44").Value = "'-4.30%"
$ws.Range("G44").Value = "'3"
$ws.Range("D45").Value = "'0.3098"
$ws.Range("E45").Value = "'0.93%"
$ws.Range("G45").Value = "'3"
$ws.Range("D46").Value = "'0.00006604"
$ws.Range("E46").Value = "'-0.98%"
$ws.Range("G46").Value = "'3"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.01%"
$ws.Range("G47").Value = "'3"
$ws.Range("E48").Value = "'31.83%"
$ws.Range("G48").Value = "'3"
$ws.Range("D49").Value = "'0.004207"
$ws.Range("E49").Value = "'2.61%"
$ws.Range("G49").Value = "'3"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("G50").Value = "'3"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'-0.01%"
$ws.Range("G51").Value = "'3"
